$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$d = $p.Designs.Item(1)
try {
  $s1.Design = $d
  Write-Output "OK-assign-design"
} catch { Write-Output ("ERR:" + $_.Exception.Message) }
